$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")
foreach ($col in $cols) {
    $r2 = $col + "2"
    $r3 = $col + "3"
    $v2 = $ws.Range($r2).Value()
    $v3 = $ws.Range($r3).Value()
    $ws.Range($r2).Value = $v3
    $ws.Range($r3).Value = $v2
}

# Swap values between row 4 and row 5 for columns D, J, K, L, M, P
foreach ($col in $cols) {
    $r4 = $col + "4"
    $r5 = $col + "5"
    $v4 = $ws.Range($r4).Value()
    $v5 = $ws.Range($r5).Value()
    $ws.Range($r4).Value = $v5
    $ws.Range($r5).Value = $v4
}
